# Iran weekly-deaths prediction workbook update (r = 9)
# Adds a new "prediction day" block for 2021-01-02 (9 rows of weekly
# forecasts), inserted between the existing 2020-12-26 and 2021-01-09
# blocks, and fills in the now-known "Real" value / error metrics for
# the first row of that new block (the week that has since completed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new blank rows at row 50, pushing the existing
# "2021-01-09" block (old rows 50-58) down to rows 59-67.
$ws.Range("A50:K58").EntireRow.Insert()

# The "day the prediction is made" column holds the date as plain text
# (e.g. "2021-01-02"), not a real Excel date. Force text format first
# so Excel doesn't auto-convert the assigned string into a date value.
$ws.Range("A50:A58").NumberFormat = "@"

# ----- New block: prediction day 2021-01-02 (rows 50-58) -----

# Row 50: week "03 Jan -- 09 Jan 2021" - now has a known Real value and
# computed error-metric columns, like other completed weeks.
$ws.Range("A50").Value2 = "2021-01-02"
$ws.Range("B50").Value2 = "03 Jan -- 09 Jan 2021"
$ws.Range("C50").Value2 = 94.56999999999999
$ws.Range("D50").Value2 = 236.37
$ws.Range("E50").Value2 = 141.8
$ws.Range("F50").Value2 = "KNN"
$ws.Range("G50").Value2 = 6.03
$ws.Range("H50").Value2 = 48.15
$ws.Range("I50").Value2 = 59.43
$ws.Range("J50").Value2 = 94.26000000000001
$ws.Range("K50").Value2 = 95.58

# Rows 51-58: remaining future weeks for the 2021-01-02 prediction,
# same structure/values as the following block (2021-01-09), since the
# underlying model forecast didn't change for those weeks.
$weeks = @(
    @{Row=51; Week="10 Jan -- 16 Jan 2021"; D=264.28},
    @{Row=52; Week="17 Jan -- 23 Jan 2021"; D=294.34},
    @{Row=53; Week="24 Jan -- 30 Jan 2021"; D=295.89},
    @{Row=54; Week="31 Jan -- 06 Feb 2021"; D=288.23},
    @{Row=55; Week="07 Feb -- 13 Feb 2021"; D=273.6},
    @{Row=56; Week="14 Feb -- 20 Feb 2021"; D=277.5},
    @{Row=57; Week="21 Feb -- 27 Feb 2021"; D=257.79},
    @{Row=58; Week="28 Feb -- 06 Mar 2021"; D=253.42}
)

foreach ($w in $weeks) {
    $r = $w.Row
    $ws.Range("A$r").Value2 = "2021-01-02"
    $ws.Range("B$r").Value2 = $w.Week
    $ws.Range("D$r").Value2 = $w.D
    $ws.Range("F$r").Value2 = "KNN"
}

$wb.Save()
